$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 16:40"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1749160
$ws.Range("C4").Value = 3357
$ws.Range("E4").Value = 1156663
$ws.Range("G4").Value = 134
$ws.Range("H4").Value = 102241

# --- Row 32: Suiza ---
$ws.Range("E32").Value = 577
$ws.Range("G32").Value = 2
$ws.Range("H32").Value = 1919

# --- Rows 62-64: Moldavia moves above Marruecos & Malasia in the ranking ---
# Row 62 becomes Moldavia with fresh data
$ws.Range("A62").Value = "Moldavia"
$ws.Range("B62").Value = 7725
$ws.Range("C62").Value = 188
$ws.Range("D62").Value = 4123
$ws.Range("E62").Value = 3320
$ws.Range("G62").Value = 8
$ws.Range("H62").Value = 282

# Row 63 becomes Marruecos (the old row-62 figures)
$ws.Range("A63").Value = "Marruecos"
$ws.Range("B63").Value = 7636
$ws.Range("C63").Value = 35
$ws.Range("D63").Value = 5109
$ws.Range("E63").Value = 2325
$ws.Range("H63").Value = 202

# Row 64 becomes Malasia (the old row-63 figures)
$ws.Range("A64").Value = "Malasia"
$ws.Range("B64").Value = 7629
$ws.Range("C64").Value = 10
$ws.Range("D64").Value = 6169
$ws.Range("E64").Value = 1345
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 115

# --- Row 77: Uzbekistan ---
$ws.Range("D77").Value = 2694
$ws.Range("E77").Value = 729

# --- Row 218: Lesoto ---
$ws.Range("D218").Value = 1
$ws.Range("E218").Value = 1
